$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 3 ----
$ws.Range("A3").Value = 131064793
$ws.Range("B3").Value = 91828
$ws.Range("E3").Value = 5432
$ws.Range("F3").Value = "Granticka"
$ws.Range("G3").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H3").Value = ""
$ws.Range("Q3").Value = 445757
$ws.Range("R3").Value = 7037091
$ws.Range("AC3").Value = ""

# ---- Row 4 ----
$ws.Range("A4").Value = 131064788
$ws.Range("M4").Value = ""
$ws.Range("Q4").Value = 445736
$ws.Range("R4").Value = 7037107
$ws.Range("AC4").Value = "Ringhack äldre"
$ws.Range("AE4").Value = $false

# ---- Row 5 ----
$ws.Range("A5").Value = 131064790
$ws.Range("B5").Value = 80348
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("Q5").Value = 445585
$ws.Range("R5").Value = 7037556

# ---- Row 6 ----
$ws.Range("A6").Value = 131064794
$ws.Range("B6").Value = 91828
$ws.Range("E6").Value = 5432
$ws.Range("F6").Value = "Granticka"
$ws.Range("G6").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H6").Value = ""
$ws.Range("Q6").Value = 445594
$ws.Range("R6").Value = 7037553
$ws.Range("AC6").Value = ""

# ---- Row 7 ----
$ws.Range("A7").Value = 131064789
$ws.Range("B7").Value = 57884
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1"
$ws.Range("I7").Style = "Normal"
$ws.Range("M7").Value = "födosökande"
$ws.Range("N7").Value = "observerad"
$ws.Range("Q7").Value = 445611
$ws.Range("R7").Value = 7037776

# ---- Row 8 ----
$ws.Range("A8").Value = 131064787
$ws.Range("B8").Value = 57884
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("M8").Value = "gammalt bo"
$ws.Range("Q8").Value = 445713
$ws.Range("R8").Value = 7037170
$ws.Range("AC8").Value = "Bohål? Ca 4,5m upp i grantickerötad granhögstubbe"
$ws.Range("AE8").Value = $true

Write-Host "Edits applied"
